$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table (GitHub Actions data pull).
# Column D cells whose new value looks numeric are pre-formatted as Text
# ("@") so Excel stores the exact original string (keeping things like
# trailing zeros / multiple thousands separators) instead of silently
# re-interpreting them as a number.

$ws.Range("D2").Value = "62.092.73"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "3.420.41"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.53"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.83"
$ws.Range("E6").Value = "  +4.98%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.09"
$ws.Range("E9").Value = "  +4.57%  "
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").Value = "4.007.62"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.66"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.436.14"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "62.106.75"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.41"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.49"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.89"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "3.563.62"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.29"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").Value = "  +4.73%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.71"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "30.96"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("D40").Value = "3.457.66"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.78"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.780"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").Value = "2.561.67"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.81"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("E51").Value = "  +0.03%  "
